# Powerpoint writer: consolidate text run nodes.
# Merge adjacent same-formatted runs (word + following space) into a
# single run per segment, reducing the number of <a:r> nodes while
# keeping the overall visible text and per-run formatting unchanged.

$p = $ppt.ActivePresentation

# Slide 1, Title "Header with inline code":
#   "Header" + " " + "with" + " " + "inline code"(Consolas)
# -> "Header " + "with " + "inline code"(Consolas)
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(1, 7).Text = "Header "
$tr1.Characters(8, 5).Text = "with "

# Slide 2, Title "Syntax highlighting":
#   "Syntax" + " " + "highlighting"
# -> "Syntax " + "highlighting"
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$tr2 = $sh2.TextFrame.TextRange
$tr2.Characters(1, 7).Text = "Syntax "

# Slide 3, Title "Two column slide":
#   "Two" + " " + "column" + " " + "slide"
# -> "Two " + "column " + "slide"
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Characters(1, 4).Text = "Two "
$tr3.Characters(5, 7).Text = "column "
